# Update countries & provincias Spain
# Refreshes the COVID "paises" dashboard: new timestamp, updated totals for
# a handful of countries, and Sudafrica's row re-inserted in its new sorted
# position (the table is kept sorted descending by "Casos totales").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 21:35"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Range("B4").Value = 1660183
$ws.Range("C4").Value = 15089
$ws.Range("D4").Value = 439083
$ws.Range("E4").Value = 1122688
$ws.Range("G4").Value = 765
$ws.Range("H4").Value = 98412

# --- Alemania (row 11) -------------------------------------------------------
$ws.Range("B11").Value = 179945
$ws.Range("C11").Value = 232
$ws.Range("E11").Value = 11682
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 8363

# --- Congo (row 136) ---------------------------------------------------------
$ws.Range("B136").Value = 487
$ws.Range("C136").Value = 18
$ws.Range("D136").Value = 147
$ws.Range("E136").Value = 324

# --- Sudafrica: updated totals move it up the ranking, ahead of Polonia ----
# Insert a fresh row right before Polonia (currently row 35) and fill it in
# with Sudafrica's refreshed numbers, then drop its old row further down
# (now shifted to row 39 after the insert).
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = "Sudafrica"
$ws.Range("B35").Value = 21343
$ws.Range("C35").Value = 1218
$ws.Range("D35").Value = 10104
$ws.Range("E35").Value = 10832
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 10
$ws.Range("H35").Value = 407
$ws.Rows.Item(39).Delete()
